$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 896, pushing the existing rows 896-937 down to
# 897-938 (dimension grows from A1:D937 to A1:D938). The new row holds
# (2026/02/27, 金, 16, 201), which sits before the row that used to be 896
# (2026/12/29, 火, 13, 88).
$ws.Rows.Item(896).Insert()

# Column A stores these dates as literal text (e.g. "2026/02/27"), not real
# date values. Assigning the string straight to .Value would get
# auto-coerced into a date serial number, so instead copy the identical
# text already sitting in A895 ("2026/02/27") down into A896 - this carries
# the text representation over without picking up any extra number-format
# style, matching the plain (unstyled) data cells elsewhere in the column.
$ws.Range("A895").Copy($ws.Range("A896"))

$ws.Cells.Item(896, 2).Value = "金"
$ws.Cells.Item(896, 3).Value = 16
$ws.Cells.Item(896, 4).Value = 201

Write-Output "inserted row 896"
